# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (Home team "H") target depth data update
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 195
$wsOff.Range("C2").Value = 144
$wsOff.Range("D2").Value = 53
$wsOff.Range("E2").Value = 24

# Sheet "DEF" - row 2 (Home team "H") target depth data update
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 259
$wsDef.Range("C2").Value = 183
$wsDef.Range("D2").Value = 44
$wsDef.Range("E2").Value = 20
